$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the full data table (header + 20 data rows) to match the updated
# public exposure sites list. New rows are added for Cheltenham, Lakes
# Entrance (x4), Mentone/Bunnings (x2), Mount Waverley and Oakleigh/Bunnings,
# and the Mordialloc 28/12 time is normalised from "2pm" to "2:00pm".

# Row 1
$ws.Range("A1").Value = 'Location'
$ws.Range("B1").Value = 'Site'
$ws.Range("C1").Value = 'Exposure period'
$ws.Range("D1").Value = 'Notes'

# Row 2
$ws.Range("A2").Value = 'Cape Schank'
$ws.Range("B2").Value = 'National Golf Club  The Cups Drive, Cape Schanck VIC 3939'
$ws.Range("C2").Value = '30/12/20 11.40am-1.40pm'
$ws.Range("D2").Value = 'Case attended course'

# Row 3
$ws.Range("A3").Value = 'Cheltenham'
$ws.Range("B3").Value = 'Aldi Cheltenham  280/282 Bay Road, Cheltenham VIC 3192'
$ws.Range("C3").Value = '29/12/2020 01:30pm-01:45pm'
$ws.Range("D3").Value = 'Case shopped in store'

# Row 4
$ws.Range("A4").Value = 'Forest Hill'
$ws.Range("B4").Value = 'Forest Hill Chase Shopping Centre 270 Canterbury Rd, Forest Hill VIC 3131'
$ws.Range("C4").Value = '28/12/20 12:00pm-2:00pm'
$ws.Range("D4").Value = '1210hrs Food court 30min; 1250hrs TKMaxx 15min; 1310hrs Target 20min; 1340hrs Woolworths 15min'

# Row 5
$ws.Range("A5").Value = 'Fountain Gate Shopping Centre'
$ws.Range("B5").Value = 'Kmart, Big W, Target, Millers, King of Gifts, Lo Costa  25-55 Overland Drive, Narre Warren VIC 3805'
$ws.Range("C5").Value = '26/12/20 9:00am-11:00am'
$ws.Range("D5").Value = ""

# Row 6
$ws.Range("A6").Value = 'Glen Waverley'
$ws.Range("B6").Value = 'Mocha Jos  87 Kingsway, Glen Waverley VIC 3150'
$ws.Range("C6").Value = '28/12/20 1:30pm-1:45pm'
$ws.Range("D6").Value = ""

# Row 7
$ws.Range("A7").Value = 'Hallam'
$ws.Range("B7").Value = 'Coles Hallam  2 Princes Domain Drive, Hallam VIC 3803'
$ws.Range("C7").Value = '30/12/20 6:15am-6:30am'
$ws.Range("D7").Value = 'Case shopped in store'

# Row 8
$ws.Range("A8").Value = 'Lakes Entrance'
$ws.Range("B8").Value = 'Blue Riviera Hire Boats  Marine Parade, Lakes Entrance VIC 3909'
$ws.Range("C8").Value = '29/12/2020 11:15am-12:15pm'
$ws.Range("D8").Value = 'Case hired a boat'

# Row 9
$ws.Range("A9").Value = 'Lakes Entrance'
$ws.Range("B9").Value = 'Central Hotel Lakes Entrance  321 Esplanade, Lakes Entrance VIC 3909'
$ws.Range("C9").Value = '30/12/2020 5:00pm-6:30pm'
$ws.Range("D9").Value = 'Case attended outside premises'

# Row 10
$ws.Range("A10").Value = 'Lakes Entrance'
$ws.Range("B10").Value = 'Darcey Annas Beach Cafe Kiosk Gift Shop Gallery  426 Main Beach Walk Surf Life Saving, Lakes Entrance VIC 3909'
$ws.Range("C10").Value = '30/12/2020 11:15am-11:20am'
$ws.Range("D10").Value = 'Case picked up takeaway'

# Row 11
$ws.Range("A11").Value = 'Lakes Entrance'
$ws.Range("B11").Value = 'Woolworths Lakes Entrance 371 Esplanade, Lakes Entrance VIC 3909'
$ws.Range("C11").Value = '30/12/2020 6:00pm-6:15pm'
$ws.Range("D11").Value = 'Case shopped in store'

# Row 12
$ws.Range("A12").Value = 'Mentone'
$ws.Range("B12").Value = 'Bunnings Mentone  23-27 Nepean Hwy, Mentone VIC 3194'
$ws.Range("C12").Value = '29/12/2020 07:30am-08:00am'
$ws.Range("D12").Value = 'Case shopped in store'

# Row 13
$ws.Range("A13").Value = 'Mentone'
$ws.Range("B13").Value = 'Bunnings Mentone  23-27 Nepean Hwy, Mentone VIC 3194'
$ws.Range("C13").Value = '31/12/2020 08:00am-08:30am'
$ws.Range("D13").Value = 'Case shopped in store'

# Row 14
$ws.Range("A14").Value = 'Mentone'
$ws.Range("B14").Value = 'Mentone/Parkdale Beach'
$ws.Range("C14").Value = '27/12/20 10:00am-4:30pm'
$ws.Range("D14").Value = ""

# Row 15
$ws.Range("A15").Value = 'Moorabbin'
$ws.Range("B15").Value = 'COSTCO Moorabbin  8 Chifley Drive, Moorabbin Airport VIC 3194'
$ws.Range("C15").Value = '30/12/20 10:45am- 12:15pm'
$ws.Range("D15").Value = 'Case shopped in store'

# Row 16
$ws.Range("A16").Value = 'Mordialloc'
$ws.Range("B16").Value = 'Woodlands Golf Club  109 White Street Mordialloc VIC 3195'
$ws.Range("C16").Value = '23/12/20 8:00am-2:00pm'
$ws.Range("D16").Value = 'Case attended course'

# Row 17
$ws.Range("A17").Value = 'Mordialloc'
$ws.Range("B17").Value = 'Woodlands Golf Club  109 White Street Mordialloc VIC 3195'
$ws.Range("C17").Value = '28/12/20 2:00pm-5:00pm'
$ws.Range("D17").Value = 'Case attended course'

# Row 18
$ws.Range("A18").Value = 'Mount Waverley'
$ws.Range("B18").Value = 'IGA  283 Stephensons Road, Mount Waverley VIC 3149'
$ws.Range("C18").Value = '21/12/2020 2:00pm-2:30pm'
$ws.Range("D18").Value = 'Case shopped for half an hour'

# Row 19
$ws.Range("A19").Value = 'Oakleigh'
$ws.Range("B19").Value = 'Bunnings Oakleigh  1041 Centre Road, Oakleigh South'
$ws.Range("C19").Value = '30/12/2020 11:00am-11:30am'
$ws.Range("D19").Value = 'Case shopped for 30 minutes'

# Row 20
$ws.Range("A20").Value = 'Oakleigh'
$ws.Range("B20").Value = 'Katialo restaurant  8 Eaton Mall, Oakleigh VIC 3166'
$ws.Range("C20").Value = '28/12/20 7:00pm-7:10pm'
$ws.Range("D20").Value = ""

# Row 21
$ws.Range("A21").Value = 'Wonthaggi'
$ws.Range("B21").Value = 'Wonthaggi Plaza Shopping centre  2 Biggs Drive, Wonthaggi VIC 3995'
$ws.Range("C21").Value = '28/12/20 1:30pm-2.30pm'
$ws.Range("D21").Value = 'Kmart- shopped for 15 mins'

